# expenses.xlsx — add "Январь 2023 г." (January 2023) as a new trailing
# month row, and correct the previous month's (row 133, "Декабрь 2022 г.")
# figures now that the real numbers for that month are known.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix up the existing last data row (Декабрь 2022 г., row 133) ---
$ws.Range("B133").Value = 151885.5
$ws.Range("C133").Value = 11946.401

# --- 2. Append a new row (134) for Январь 2023 г. -----------------------
# Copy the formatting (styles + row height) of the row above down into the
# new row before writing its values, so the new row matches the rest of
# the table.
$ws.Range("A133:C133").Copy()
$ws.Range("A134:C134").PasteSpecial(-4122)
$ws.Range("A134:C134").RowHeight = 11.45
$excel.CutCopyMode = 0

$ws.Range("A134").Value = "Январь 2023 г."
$ws.Range("B134").Value = 66284.399999999994
$ws.Range("C134").Value = 7809.7

# --- 3. Match the author's final selection (cosmetic) --------------------
$ws.Range("H104").Select() | Out-Null
